$d = $word.ActiveDocument

# The ID-marker paragraph is the first paragraph in the document.
$para = $d.Paragraphs(1)

# New left indent: 225 twips == 11.25 pt (Word's ParagraphFormat properties
# are expressed in points, which get persisted as twentieths-of-a-point).
$para.Range.ParagraphFormat.LeftIndent = 11.25

# Add a paragraph border box that reserves 5pt of padding on every side but
# draws no visible rule (LineStyle = wdLineStyleNone).
$borders = $para.Range.ParagraphFormat.Borders
$borders.Item(-1).LineStyle = 0
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Replace the ID placeholder text ("...topic_3..." -> "...202...") and fold
# the trailing stand-alone space run into the same run by replacing the
# whole "text + space" span in one shot.
$d.Content.Find.Execute("**ID__AFFARS_5306_topic_3__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5306_202__ID**", 2)
